$d = $word.ActiveDocument

# Locate the paragraph containing "Mambo Marie" (the heading line), then
# the empty paragraph immediately after it is where the first two new
# paragraphs of body text get inserted (replacing the lone empty <w:p/>).
$marieRange = $d.Content
$found = $marieRange.Find.Execute("Mambo Marie", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mariePara = $marieRange.Paragraphs(1)
$targetPara = $mariePara.Next()
$targetRange = $targetPara.Range

$block1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve"> The Mambo Marie feature was implemented in</w:t></w:r><w:r><w:t xml:space="preserve">to the game. To represent a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MamboMarie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> object, a class named “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MamboMarie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>” was created. This class extends from the “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ZombieActor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">” class. By doing this, we can easily access and use critical methods that are already defined for us in the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ZombieActor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class, and therefore also the Actor class. An alternative to this design decision may be to have the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MamboMarie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class inherit Zombie rather than </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ZombieActor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">This was not done as the Zombie class does not have a constructor (by default or from previous assignments) where we can choose the display character, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hitpoints</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and other attributes</w:t></w:r><w:r><w:t xml:space="preserve">. We could overload the constructor in Zombie which will let us do this, however I decided that the Zombie class did not have any methods that would be useful to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>MamboMarie</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ZombieActor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> already had this sort of constructor defined.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetRange.InsertXML($block1)

# Locate the bookmarked paragraph (_GoBack) -- it now sits right after the
# two paragraphs we just inserted -- and add two blank paragraphs right
# after it, before the "Ending the game" heading.
$bm = $d.Bookmarks("_GoBack")
$bmPara = $bm.Range.Paragraphs(1)
$insertPos = $bmPara.Range.End
$gapRange = $d.Range($insertPos, $insertPos)

$block2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$gapRange.InsertXML($block2)
